$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "Color" column (currently column U) and
# give the new header cell the same style as the other plain headers (s="4"),
# which the column-insert operation already inherits from the column to its
# right.
$ws.Range("U1").EntireColumn.Insert()
$ws.Range("U1").Value = "thickness"

# Match the new selection left behind by the edit (Excel recorded the cursor
# at U4 when the file was re-saved).
$ws.Range("U4").Select()
